$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 16 corresponds to group "deloitte_exam_answers"
# Update current_phase (D16): 1 -> 2
$ws.Range("D16").Value = 2

# Update last_action_date (E16): empty -> timestamp
$ws.Range("E16").Value = "2026-02-12T10:15:28.129547+00:00"

# Update reactions_count (H16): 0 -> 2
$ws.Range("H16").Value = 2

# Update replies_count (I16): 0 -> 1
$ws.Range("I16").Value = 1

# Update reacted_message_ids (L16): [] -> [169, 425]
$ws.Range("L16").Value = "[169, 425]"

# Update replied_message_ids (M16): [] -> [148]
$ws.Range("M16").Value = "[148]"
